$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the new "table_definitions" sheet as the very first tab.
# ---------------------------------------------------------------------------
$firstSheet = $wb.Worksheets.Item(1)
$ws = $wb.Worksheets.Add($firstSheet)
$ws.Name = "table_definitions"

# ---------------------------------------------------------------------------
# 2. Populate the sheet, row by row, left-to-right within each row, in the
#    same order the values were originally typed in (rows 1-5, then 8, then
#    6-7) so that the shared-strings table comes out in the same order.
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "mapping_file_name"
$ws.Range("B1").Value = "entity_name"
$ws.Range("C1").Value = "required_entities"
$ws.Range("D1").Value = "destination_table_name"
$ws.Range("E1").Value = "table_type"
$ws.Range("F1").Value = "source_table_name"
$ws.Range("G1").Value = "casrec_conditions"
$ws.Range("H1").Value = "source_table_additional_columns"

$ws.Range("A2").Value = "client_nodebtchase_warnings"
$ws.Range("B2").Value = "warnings"
$ws.Range("D2").Value = "warnings"
$ws.Range("E2").Value = "data"
$ws.Range("F2").Value = "pat"
$ws.Range("G2").Value = "Debt chase = not null"
$ws.Range("H2").Value = "Case"

$ws.Range("A3").Value = "client_saarcheck_warnings"
$ws.Range("B3").Value = "warnings"
$ws.Range("D3").Value = "warnings"
$ws.Range("E3").Value = "data"
$ws.Range("F3").Value = "pat"
$ws.Range("G3").Value = "SAAR Check = not null"
$ws.Range("H3").Value = "Case"

$ws.Range("A4").Value = "client_special_warnings"
$ws.Range("B4").Value = "warnings"
$ws.Range("D4").Value = "warnings"
$ws.Range("E4").Value = "data"
$ws.Range("F4").Value = "pat"
$ws.Range("G4").Value = "SIM = not null"
$ws.Range("H4").Value = "Case"

$ws.Range("A5").Value = "client_violent_warnings"
$ws.Range("B5").Value = "warnings"
$ws.Range("D5").Value = "warnings"
$ws.Range("E5").Value = "data"
$ws.Range("F5").Value = "pat"
$ws.Range("G5").Value = "VWM = not null"
$ws.Range("H5").Value = "Case"

$ws.Range("A8").Value = "person_warning"
$ws.Range("B8").Value = "warnings"
$ws.Range("D8").Value = "person_warning"
$ws.Range("E8").Value = "join"

$ws.Range("A6").Value = "deputy_special_warnings"
$ws.Range("B6").Value = "warnings"
$ws.Range("D6").Value = "warnings"
$ws.Range("E6").Value = "data"
$ws.Range("F6").Value = "deputy"
$ws.Range("G6").Value = "SIM = not null"
$ws.Range("H6").Value = "Deputy No"

$ws.Range("A7").Value = "deputy_violent_warnings"
$ws.Range("B7").Value = "warnings"
$ws.Range("D7").Value = "warnings"
$ws.Range("E7").Value = "data"
$ws.Range("F7").Value = "deputy"
$ws.Range("G7").Value = "VWM = not null"
$ws.Range("H7").Value = "Deputy No"

# ---------------------------------------------------------------------------
# 3. Formatting - dark-grey Helvetica header row, black Arial body rows (both
#    size 10). Multi-area ranges are used so cells that were never given a
#    value (column C throughout, and F8:H8) are not touched.
# ---------------------------------------------------------------------------
$headerFont = $ws.Range("A1:H1").Font
$headerFont.Name = "Helvetica"
$headerFont.Size = 10
$headerFont.Color = 3355443

$bodyRange = $ws.Range("A2:B2,D2:H2,A3:B3,D3:H3,A4:B4,D4:H4,A5:B5,D5:H5,A6:B6,D6:H6,A7:B7,D7:H7,A8:B8,D8:E8")
$bodyFont = $bodyRange.Font
$bodyFont.Name = "Arial"
$bodyFont.Size = 10
$bodyFont.Color = 0

# ---------------------------------------------------------------------------
# 4. Sheet view - zoomed in, tab selected/active, cursor on G6.
# ---------------------------------------------------------------------------
$null = $ws.Activate()
$zoomWindow = $excel.ActiveWindow
$zoomWindow.Zoom = 200
$null = $ws.Range("G6").Select()
